# Updates scraped market-price data across the Chocobo_Profits workbook.
# Each leve's currentAveragePrice / NQ / HQ price & profit columns (H:N)
# are refreshed with new values from the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9918.888999999999
$ws.Range("I51").Value = 8300
$ws.Range("J51").Value = 10121.25
$ws.Range("K51").Value = 8300
$ws.Range("L51").Value = 10121.25
$ws.Range("M51").Value = -7816
$ws.Range("N51").Value = -11089.25
$ws.Range("H53").Value = 550.35
$ws.Range("I53").Value = 265.72726
$ws.Range("J53").Value = 898.2222
$ws.Range("K53").Value = 265.72726
$ws.Range("L53").Value = 898.2222
$ws.Range("M53").Value = 371.27274
$ws.Range("N53").Value = -2172.2222
$ws.Range("H101").Value = 3227
$ws.Range("I101").Value = 459
$ws.Range("J101").Value = 5995
$ws.Range("K101").Value = 1377
$ws.Range("L101").Value = 17985
$ws.Range("M101").Value = 245
$ws.Range("N101").Value = -21229
$ws.Range("H138").Value = 2449.0505
$ws.Range("I138").Value = 669.17645
$ws.Range("J138").Value = 2818.0488
$ws.Range("K138").Value = 2007.52935
$ws.Range("L138").Value = 8454.1464
$ws.Range("M138").Value = 3132.47065
$ws.Range("N138").Value = -18734.1464

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 571.4286
$ws.Range("I2").Value = 571.4286
$ws.Range("K2").Value = 571.4286
$ws.Range("M2").Value = -458.4286
$ws.Range("H32").Value = 3596.7605
$ws.Range("I32").Value = 2972.111
$ws.Range("K32").Value = 2972.111
$ws.Range("M32").Value = -2685.111
$ws.Range("H74").Value = 8114.357
$ws.Range("I74").Value = 9923
$ws.Range("J74").Value = 4858.8
$ws.Range("K74").Value = 9923
$ws.Range("L74").Value = 4858.8
$ws.Range("M74").Value = -9049
$ws.Range("N74").Value = -6606.8
$ws.Range("H77").Value = 8114.357
$ws.Range("I77").Value = 9923
$ws.Range("J77").Value = 4858.8
$ws.Range("K77").Value = 49615
$ws.Range("L77").Value = 24294
$ws.Range("M77").Value = -45247
$ws.Range("N77").Value = -33030
$ws.Range("H97").Value = 1403.45
$ws.Range("I97").Value = 947.6667
$ws.Range("J97").Value = 5505.5
$ws.Range("K97").Value = 947.6667
$ws.Range("L97").Value = 5505.5
$ws.Range("M97").Value = -451.6667
$ws.Range("N97").Value = -6497.5
$ws.Range("H116").Value = 571.4286
$ws.Range("I116").Value = 571.4286
$ws.Range("K116").Value = 571.4286
$ws.Range("M116").Value = 1722.5714
$ws.Range("H132").Value = 2125.4443
$ws.Range("I132").Value = 939.86664
$ws.Range("J132").Value = 3607.4167
$ws.Range("K132").Value = 2819.59992
$ws.Range("L132").Value = 10822.2501
$ws.Range("M132").Value = -289.5999199999997
$ws.Range("N132").Value = -15882.2501

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 571.4286
$ws.Range("I3").Value = 571.4286
$ws.Range("K3").Value = 571.4286
$ws.Range("M3").Value = -457.4286
$ws.Range("H54").Value = 3279.4
$ws.Range("I54").Value = 3279.4
$ws.Range("K54").Value = 3279.4
$ws.Range("M54").Value = -2795.4
$ws.Range("H105").Value = 1655.0139
$ws.Range("I105").Value = 1631.8383
$ws.Range("J105").Value = 2049
$ws.Range("K105").Value = 1631.8383
$ws.Range("L105").Value = 2049
$ws.Range("M105").Value = 115.1617000000001
$ws.Range("N105").Value = -5543

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1390.3928
$ws.Range("I105").Value = 1117.24
$ws.Range("K105").Value = 1117.24
$ws.Range("M105").Value = 629.76
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H132").Value = 2339.423
$ws.Range("I132").Value = 1438.2632
$ws.Range("J132").Value = 4785.4287
$ws.Range("K132").Value = 4314.7896
$ws.Range("L132").Value = 14356.2861
$ws.Range("M132").Value = -1784.7896
$ws.Range("N132").Value = -19416.2861
$ws.Range("H134").Value = 7094.05
$ws.Range("I134").Value = 9373.583000000001
$ws.Range("K134").Value = 28120.749
$ws.Range("M134").Value = -25585.749
$ws.Range("N110").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 607752.5600000001
$ws.Range("I5").Value = 651.2857
$ws.Range("J5").Value = 891066.5600000001
$ws.Range("K5").Value = 1953.8571
$ws.Range("L5").Value = 2673199.68
$ws.Range("M5").Value = -1841.8571
$ws.Range("N5").Value = -2673423.68
$ws.Range("H92").Value = 823.625
$ws.Range("I92").Value = 827
$ws.Range("J92").Value = 800
$ws.Range("K92").Value = 2481
$ws.Range("L92").Value = 2400
$ws.Range("M92").Value = -1233
$ws.Range("N92").Value = -4896
$ws.Range("H105").Value = 4999
$ws.Range("J105").Value = 4999
$ws.Range("L105").Value = 14997
$ws.Range("N105").Value = -20239
$ws.Range("H113").Value = 3906833.5
$ws.Range("I113").Value = 618.4286
$ws.Range("J113").Value = 6945000.5
$ws.Range("K113").Value = 1855.2858
$ws.Range("L113").Value = 20835001.5
$ws.Range("M113").Value = 314.7142000000001
$ws.Range("N113").Value = -20839341.5
$ws.Range("H135").Value = 607752.5600000001
$ws.Range("I135").Value = 651.2857
$ws.Range("J135").Value = 891066.5600000001
$ws.Range("K135").Value = 5861.571300000001
$ws.Range("L135").Value = 8019599.040000001
$ws.Range("M135").Value = -3326.571300000001
$ws.Range("N135").Value = -8024669.040000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 35666.668
$ws.Range("J120").Value = 35666.668
$ws.Range("L120").Value = 35666.668
$ws.Range("N120").Value = -45342.668
$ws.Range("H132").Value = 2987.7307
$ws.Range("I132").Value = 1643.7858
$ws.Range("J132").Value = 4555.6665
$ws.Range("K132").Value = 4931.357400000001
$ws.Range("L132").Value = 13666.9995
$ws.Range("M132").Value = -2401.357400000001
$ws.Range("N132").Value = -18726.9995

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9261601
$ws.Range("I93").Value = 15874931
$ws.Range("J93").Value = 2940
$ws.Range("K93").Value = 15874931
$ws.Range("L93").Value = 2940
$ws.Range("M93").Value = -15873683
$ws.Range("N93").Value = -5436
$ws.Range("H132").Value = 5334.4
$ws.Range("I132").Value = 992.2857
$ws.Range("J132").Value = 15466
$ws.Range("K132").Value = 2976.8571
$ws.Range("L132").Value = 46398
$ws.Range("M132").Value = -446.8571000000002
$ws.Range("N132").Value = -51458

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("H132").Value = 9011822
$ws.Range("I132").Value = 1813.4762
$ws.Range("J132").Value = 20837458
$ws.Range("K132").Value = 5440.4286
$ws.Range("L132").Value = 62512374
$ws.Range("M132").Value = -2910.4286
$ws.Range("N132").Value = -62517434
$ws.Range("M61").ClearContents()
